$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column retains its original text representation
# (values like "1.001" or "20.50" would otherwise be auto-converted to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.976.07'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '1.782.99'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '316.04'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '0.5372'
$ws.Range("E7").Value = '  -2.09%  '
$ws.Range("D8").Value = '0.3756'
$ws.Range("E8").Value = '  -2.54%  '
$ws.Range("D9").Value = '0.07442'
$ws.Range("E9").Value = '  -2.10%  '
$ws.Range("D10").Value = '41.69'
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("D11").Value = '1.090'
$ws.Range("E11").Value = '  -3.23%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '20.50'
$ws.Range("E13").Value = '  -3.14%  '
$ws.Range("D14").Value = '6.086'
$ws.Range("E14").Value = '  -1.79%  '
$ws.Range("D15").Value = '7.205'
$ws.Range("E15").Value = '  -2.44%  '
$ws.Range("D16").Value = '1.776.57'
$ws.Range("E16").Value = '  -1.41%  '
$ws.Range("D17").Value = '88.55'
$ws.Range("E17").Value = '  -4.04%  '
$ws.Range("D18").Value = '0.00001054'
$ws.Range("E18").Value = '  -1.62%  '
$ws.Range("D19").Value = '0.06441'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").Value = '17.26'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("D22").Value = '5.884'
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").Value = '28.002.46'
$ws.Range("E23").Value = '  -1.47%  '
$ws.Range("D24").Value = '11.22'
$ws.Range("E24").Value = '  -1.98%  '
$ws.Range("D25").Value = '2.088'
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").Value = '155.93'
$ws.Range("E26").Value = '  -1.71%  '
$ws.Range("D27").Value = '20.23'
$ws.Range("E27").Value = '  -2.11%  '
$ws.Range("D28").Value = '1.982.06'
$ws.Range("E28").Value = '  -1.49%  '
$ws.Range("D29").Value = '2.273'
$ws.Range("E29").Value = '  -5.52%  '
$ws.Range("D30").Value = '119.94'
$ws.Range("E30").Value = '  -3.10%  '
$ws.Range("D31").Value = '1.105'
$ws.Range("E31").Value = '  -1.89%  '
$ws.Range("D32").Value = '0.1050'
$ws.Range("E32").Value = '  +2.91%  '
$ws.Range("D33").Value = '3.638'
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("D34").Value = '5.527'
$ws.Range("E34").Value = '  -3.99%  '
$ws.Range("D35").Value = '0.2259'
$ws.Range("E35").Value = '  -2.89%  '
$ws.Range("D36").Value = '0.06423'
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").Value = '0.02279'
$ws.Range("E37").Value = '  -2.06%  '
$ws.Range("D38").Value = '5.002'
$ws.Range("E38").Value = '  -1.94%  '
$ws.Range("D39").Value = '8.445'
$ws.Range("E39").Value = '  -4.61%  '
$ws.Range("D40").Value = '1.444'
$ws.Range("E40").Value = '  +4.36%  '
$ws.Range("D41").Value = '0.6136'
$ws.Range("E41").Value = '  -4.66%  '
$ws.Range("D42").Value = '11.06'
$ws.Range("E42").Value = '  -5.02%  '
$ws.Range("D43").Value = '1.172'
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = '13.29'
$ws.Range("E45").Value = '  -1.55%  '
$ws.Range("D46").Value = '3.664'
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("D47").Value = '0.5742'
$ws.Range("E47").Value = '  -4.06%  '
$ws.Range("D48").Value = '126.58'
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").Value = '1.186'
$ws.Range("E49").Value = '  +3.32%  '
$ws.Range("D50").Value = '1.921'
$ws.Range("E50").Value = '  -3.63%  '
$ws.Range("D51").Value = '0.06794'
$ws.Range("E51").Value = '  -1.62%  '
